# Update the "想去人数" (number of people interested) figures that changed
# between the two site generation runs.
#
# Sheet "展览" (exhibitions):
#   F3: 2136 -> 2162
#   F4: 878  -> 885
#   F5: 1419 -> 1444
#   F6: 373  -> 375
#
# Sheet "全部类型" (all types), same events appear at different rows:
#   F3: 2136 -> 2162
#   F6: 878  -> 885
#   F7: 1419 -> 1444
#   F8: 373  -> 375

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F3").Value = 2162
$wsExhibit.Range("F4").Value = 885
$wsExhibit.Range("F5").Value = 1444
$wsExhibit.Range("F6").Value = 375

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F3").Value = 2162
$wsAll.Range("F6").Value = 885
$wsAll.Range("F7").Value = 1444
$wsAll.Range("F8").Value = 375
